$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new daily row (47) below the existing data.
# Column A holds a date-like string ("2025/10/02"); format the cell as
# Text first so Excel's automatic date-recognition doesn't convert it to
# a serial date number, then reset the style back to the sheet's default
# so the new row stays unstyled like the rows above it.
$rowNum = 47
$a = $ws.Cells.Item($rowNum, 1)
$a.NumberFormat = "@"
$a.Value = "2025/10/02"
$a.Style = "Normal"

$ws.Cells.Item($rowNum, 2).Value = "木"
$ws.Cells.Item($rowNum, 3).Value = 0
$ws.Cells.Item($rowNum, 4).Value = 3
